# Auto-generated edit script applying the diff to Ultros_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 700
$ws.Range("I7").Value = 400
$ws.Range("J7").Value = 1000
$ws.Range("K7").Value = 400
$ws.Range("L7").Value = 1000
$ws.Range("N7").Value = -1224
$ws.Range("M7").Value = -288

$ws.Range("H10").Value = 496.33334
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 496.33334
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 496.33334
$ws.Range("M10").Value = ""
$ws.Range("N10").Value = -1082.33334

$ws.Range("H12").Value = 345
$ws.Range("J12").Value = 500
$ws.Range("L12").Value = 500
$ws.Range("N12").Value = -840

$ws.Range("H14").Value = 700
$ws.Range("I14").Value = 400
$ws.Range("J14").Value = 1000
$ws.Range("K14").Value = 400
$ws.Range("L14").Value = 1000
$ws.Range("N14").Value = -1382
$ws.Range("M14").Value = -209

$ws.Range("H58").Value = 379.9
$ws.Range("J58").Value = 1250
$ws.Range("L58").Value = 3750
$ws.Range("N58").Value = -4050

$ws.Range("H70").Value = 4067.3914
$ws.Range("I70").Value = 1225
$ws.Range("K70").Value = 3675
$ws.Range("M70").Value = -3405

$ws.Range("H73").Value = 4067.3914
$ws.Range("I73").Value = 1225
$ws.Range("K73").Value = 3675
$ws.Range("M73").Value = -2739

$ws.Range("H87").Value = 19999.773
$ws.Range("J87").Value = 19999.773
$ws.Range("L87").Value = 19999.773
$ws.Range("N87").Value = -22495.773

$ws.Range("H90").Value = 19999.773
$ws.Range("J90").Value = 19999.773
$ws.Range("L90").Value = 59999.319
$ws.Range("N90").Value = -72479.319

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6898.25
$ws.Range("I32").Value = 6593.7676
$ws.Range("K32").Value = 6593.7676
$ws.Range("M32").Value = -6306.7676

$ws.Range("H45").Value = 2910.1
$ws.Range("I45").Value = 2406
$ws.Range("K45").Value = 2406
$ws.Range("M45").Value = -2029

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 291
$ws.Range("I5").Value = 237.75
$ws.Range("J5").Value = 397.5
$ws.Range("K5").Value = 237.75
$ws.Range("L5").Value = 397.5
$ws.Range("M5").Value = -125.75
$ws.Range("N5").Value = -621.5

$ws.Range("H10").Value = 638.6667
$ws.Range("I10").Value = 433
$ws.Range("K10").Value = 433
$ws.Range("M10").Value = -294

$ws.Range("H11").Value = 1050
$ws.Range("J11").Value = 1050
$ws.Range("L11").Value = 1050
$ws.Range("N11").Value = -1330

$ws.Range("H13").Value = 2120
$ws.Range("J13").Value = 3135
$ws.Range("L13").Value = 3135
$ws.Range("N13").Value = -3413

$ws.Range("H14").Value = 750
$ws.Range("J14").Value = 750
$ws.Range("L14").Value = 750
$ws.Range("N14").Value = -1090

$ws.Range("H19").Value = 176.35
$ws.Range("I19").Value = 126.5
$ws.Range("J19").Value = 625
$ws.Range("K19").Value = 126.5
$ws.Range("L19").Value = 625
$ws.Range("M19").Value = 43.5
$ws.Range("N19").Value = -965

$ws.Range("H24").Value = 176.35
$ws.Range("I24").Value = 126.5
$ws.Range("J24").Value = 625
$ws.Range("K24").Value = 126.5
$ws.Range("L24").Value = 625
$ws.Range("M24").Value = 43.5
$ws.Range("N24").Value = -965

$ws.Range("H25").Value = 2182.2
$ws.Range("I25").Value = 1502.75
$ws.Range("K25").Value = 1502.75
$ws.Range("M25").Value = -1328.75

$ws.Range("H31").Value = 5811.9165
$ws.Range("I31").Value = 1270.3334
$ws.Range("J31").Value = 7325.778
$ws.Range("K31").Value = 1270.3334
$ws.Range("L31").Value = 7325.778
$ws.Range("M31").Value = -975.3334
$ws.Range("N31").Value = -7915.778

$ws.Range("H34").Value = 5811.9165
$ws.Range("I34").Value = 1270.3334
$ws.Range("J34").Value = 7325.778
$ws.Range("K34").Value = 1270.3334
$ws.Range("L34").Value = 7325.778
$ws.Range("M34").Value = -1068.3334
$ws.Range("N34").Value = -7729.778

$ws.Range("H35").Value = 1600
$ws.Range("J35").Value = 2500
$ws.Range("L35").Value = 2500
$ws.Range("N35").Value = -3088

$ws.Range("H62").Value = 4999.1665
$ws.Range("I62").Value = 4999.1665
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 4999.1665
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -4375.1665
$ws.Range("N62").Value = ""

$ws.Range("H65").Value = 4999.1665
$ws.Range("I65").Value = 4999.1665
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 24995.8325
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -21875.8325
$ws.Range("N65").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 2325
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 2325
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 6975
$ws.Range("M25").Value = ""
$ws.Range("N25").Value = -7313

$ws.Range("H29").Value = 68519440
$ws.Range("I29").Value = 15556655
$ws.Range("J29").Value = 333333340
$ws.Range("K29").Value = 46669965
$ws.Range("L29").Value = 1000000020
$ws.Range("M29").Value = -46669688
$ws.Range("N29").Value = -1000000574

$ws.Range("H30").Value = 2325
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 2325
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 6975
$ws.Range("M30").Value = ""
$ws.Range("N30").Value = -7179

$ws.Range("H55").Value = 114444780
$ws.Range("I55").Value = 1000
$ws.Range("J55").Value = 171666670
$ws.Range("K55").Value = 3000
$ws.Range("L55").Value = 515000010
$ws.Range("M55").Value = -2823
$ws.Range("N55").Value = -515000364

$ws.Range("H81").Value = 2752.5
$ws.Range("J81").Value = 2752.5
$ws.Range("L81").Value = 8257.5
$ws.Range("N81").Value = -10503.5

$ws.Range("H84").Value = 2752.5
$ws.Range("J84").Value = 2752.5
$ws.Range("L84").Value = 24772.5
$ws.Range("N84").Value = -36004.5

$ws.Range("H104").Value = 3201.6365
$ws.Range("I104").Value = 1153.5
$ws.Range("J104").Value = 8663.333000000001
$ws.Range("K104").Value = 3460.5
$ws.Range("L104").Value = 25989.999
$ws.Range("M104").Value = -839.5
$ws.Range("N104").Value = -31231.999

$ws.Range("H121").Value = 1601.6666
$ws.Range("I121").Value = 453.33334
$ws.Range("J121").Value = 1831.3334
$ws.Range("K121").Value = 1360.00002
$ws.Range("L121").Value = 5494.0002
$ws.Range("M121").Value = -50.00001999999995
$ws.Range("N121").Value = -8114.0002

$ws.Range("H131").Value = 5399.385
$ws.Range("I131").Value = 3049.6667
$ws.Range("J131").Value = 7413.4287
$ws.Range("K131").Value = 9149.000100000001
$ws.Range("L131").Value = 22240.2861
$ws.Range("M131").Value = -4109.000100000001
$ws.Range("N131").Value = -32320.2861

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 5722.75
$ws.Range("I43").Value = 656.8
$ws.Range("J43").Value = 14166
$ws.Range("K43").Value = 656.8
$ws.Range("L43").Value = 14166
$ws.Range("M43").Value = -505.8
$ws.Range("N43").Value = -14468

$ws.Range("H46").Value = 14600
$ws.Range("J46").Value = 15000
$ws.Range("L46").Value = 15000
$ws.Range("N46").Value = -15312

$ws.Range("H63").Value = 36490.832
$ws.Range("J63").Value = 36490.832
$ws.Range("L63").Value = 36490.832
$ws.Range("N63").Value = -37862.832

$ws.Range("H66").Value = 36490.832
$ws.Range("J66").Value = 36490.832
$ws.Range("L66").Value = 109472.496
$ws.Range("N66").Value = -116336.496

$ws.Range("H113").Value = 8888.888999999999
$ws.Range("I113").Value = 4600
$ws.Range("J113").Value = 14250
$ws.Range("K113").Value = 4600
$ws.Range("L113").Value = 14250
$ws.Range("M113").Value = -2430
$ws.Range("N113").Value = -18590

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 4614.857
$ws.Range("I10").Value = 500
$ws.Range("J10").Value = 5300.6665
$ws.Range("K10").Value = 500
$ws.Range("L10").Value = 5300.6665
$ws.Range("N10").Value = -5580.6665
$ws.Range("M10").Value = -360

$ws.Range("H19").Value = 1729.1428
$ws.Range("I19").Value = 902
$ws.Range("J19").Value = 2349.5
$ws.Range("K19").Value = 902
$ws.Range("L19").Value = 2349.5
$ws.Range("M19").Value = -732
$ws.Range("N19").Value = -2689.5

$ws.Range("H46").Value = 932.6667
$ws.Range("I46").Value = 838.8
$ws.Range("J46").Value = 1050
$ws.Range("K46").Value = 838.8
$ws.Range("L46").Value = 1050
$ws.Range("M46").Value = -650.8
$ws.Range("N46").Value = -1426

$ws.Range("H81").Value = 54000
$ws.Range("J81").Value = 54000
$ws.Range("L81").Value = 54000
$ws.Range("N81").Value = -55996

$ws.Range("H84").Value = 54000
$ws.Range("J84").Value = 54000
$ws.Range("L84").Value = 162000
$ws.Range("N84").Value = -171984

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 1066.6666
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 1066.6666
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 1066.6666
$ws.Range("M8").Value = ""
$ws.Range("N8").Value = -1346.6666

$ws.Range("H56").Value = 15564.667
$ws.Range("I56").Value = 9285
$ws.Range("J56").Value = 18704.5
$ws.Range("K56").Value = 9285
$ws.Range("L56").Value = 18704.5
$ws.Range("M56").Value = -8571
$ws.Range("N56").Value = -20132.5

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").Value = ""

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").Value = ""
